$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the sheet's label cell
$ws.Name = "Tienda Velázquez"

# Update column A (Tienda) for data rows 2-11
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "Tienda Velázquez"
}

# Fill column B (Nombre_TPV) for data rows 2-11
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 2).Value = "BAR"
}
for ($r = 6; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = "SERVIDOR TIENDA"
}
